$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dosing channels (Priva point UUID in column A, friendly name in column B)
$ws.Range("A77").Value = "8000008c-0001-0001-0000-000080000a77"
$ws.Range("A78").Value = "8000008f-0001-8101-0000-0000800007d3"
$ws.Range("B78").Value = "Dosing_pH_1"
$ws.Range("B77").Value = "Dosing_EC"
$ws.Range("B79").Value = "Dosing_pH_2"
$ws.Range("A79").Value = "8000008f-0001-8201-0000-0000800007d3"

# Trailing styled-but-empty row, matching the "list of all available datapoints"
# marker style used elsewhere in the sheet (new Menlo/orange font)
$c = $ws.Range("A80")
$c.Font.Color = 7901646
$c.Font.Name = "Menlo"

# Land the selection back at the top-left, clearing the stale scrolled/selected state
$ws.Range("A1").Select()
